$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.059.26'
$ws.Range("E2").Value = '  +0.44%  '

$ws.Range("D3").Value = '1.640.00'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("E4").Value = '  +0.51%  '

$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.70'
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = '  -0.44%  '

$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("E7").Value = '  +0.45%  '

$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.252'
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = '  -1.68%  '

$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0627'
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = '  -1.55%  '

$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.69'
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = '  -4.63%  '

$ws.Range("E11").Value = '  -0.04%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.652.17'
$ws.Range("E12").Value = '  -0.01%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.21'
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = '  -1.43%  '

$ws.Range("E14").Value = '  -2.01%  '

$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '62.40'
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = '  -0.79%  '

$ws.Range("E16").Value = '  -1.89%  '

$ws.Range("D17").Value = '26.073.37'
$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("E18").Value = '  +0.51%  '

$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '191.34'
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = '  -0.68%  '

$ws.Range("E20").Value = '  -1.71%  '

$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.61'
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = '  -3.05%  '

$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.15'
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = '  -1.79%  '

$ws.Range("B23").Value = 'Stellar'
$ws.Range("C23").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.131'
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = '  -0.36%  '

$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.06'
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = '  +0.57%  '

$ws.Range("E25").Value = '  +0.56%  '

$ws.Range("E26").Value = '  -0.92%  '

$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.78'
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = '  -1.40%  '

$ws.Range("E28").Value = '  -1.93%  '

$ws.Range("E29").Value = '  -0.21%  '

$ws.Range("E30").Value = '  -2.58%  '

$ws.Range("E31").Value = '  -2.09%  '

$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.19'
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = '  -3.02%  '

$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.51'
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = '  -0.95%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.44'
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = '  -0.92%  '

$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.879'
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = '  -2.41%  '

$ws.Range("D36").Value = '1.133.40'
$ws.Range("E36").Value = '  +0.12%  '

$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.45'
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = '  -0.12%  '

$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.526'
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = '  -3.13%  '

$ws.Range("E39").Value = '  -1.30%  '

$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '98.89'
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = '  -0.36%  '

$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.786'
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = '  -1.16%  '

$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.31'
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = '  -2.95%  '

$ws.Range("D43").Value = '0.0₆0114'
$ws.Range("E43").Value = '  -0.82%  '

$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '55.46'
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = '  -2.01%  '

$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0530'
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = '  +0.19%  '

$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.50'
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = '  +1.63%  '

$ws.Range("E47").Value = '  +0.02%  '

$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.57'
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = '  -1.43%  '

$ws.Range("E49").Value = '  +0.11%  '

$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0929'
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = '  -3.20%  '

$ws.Range("E51").Value = '  -0.62%  '
